# Applies the "Add files via upload" commit:
#   - Dt. Referencia (column G) rolls forward one day for every data row
#     (45489 -> 45490, i.e. 2024-07-16 -> 2024-07-17)
#   - Three rows get corrected Saldo Previsto / Vl. Total values
#     (columns E and H)
#   - The sheet is renamed to reflect the new export timestamp

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Roll every "Dt. Referencia" (column G) value forward by one day for all
# data rows (row 2 through the last used row).
$lastRow = $ws.UsedRange.Rows.Count
$ws.Range("G2:G" + $lastRow).Value = 45490

# Row 118: Saldo Previsto / Vl. Total correction.
$ws.Range("E118").Value = 19019.49
$ws.Range("H118").Value = 19019.49

# Row 224: Saldo Previsto / Vl. Total correction.
$ws.Range("E224").Value = 476.66
$ws.Range("H224").Value = 476.66

# Row 255: Saldo Previsto / Vl. Total correction.
$ws.Range("E255").Value = 10168.82
$ws.Range("H255").Value = 10168.82

# Rename the sheet to match the new export run (2024-07-17 09:44:48).
$ws.Name = "IClientBalance-20240717-094448-"
